$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coerce cells whose new numeric-looking value would otherwise be
# auto-converted from text to a number back to Text, so the stored
# content exactly matches the source price strings (e.g. trailing
# zeros like "1.00" / "0.460" / "17.70" must be preserved as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin / price / link / volume values.
$ws.Range("D2").Value = "71.169.24"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.849.15"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "693.93"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "173.27"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "3.846.27"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "7.25"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  +4.29%  "
$ws.Range("D14").Value = "36.37"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "4.498.67"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "3.847.84"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "71.221.39"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "17.70"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "11.14"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "493.76"
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "85.03"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "12.27"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "10.57"
$ws.Range("E27").Value = "  +2.26%  "
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "4.004.52"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +8.30%  "
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "29.62"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "0.179"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").Value = "3.801.64"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("E40").Value = "  +12.54%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "6.04"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("E43").Value = "  +6.36%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D46").Value = "163.80"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "48.61"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "44.17"
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "0.303"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "417.68"
$ws.Range("E51").Value = "  +4.70%  "
